# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Row 2 holds the "property" identifiers for each column, row 3 the
# property "kind" (medida/dim), row 4 the datatype/URI template, and
# row 5 held per-column mapping-file references that are no longer
# needed now that the dimensions have been curated directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: property identifiers
$ws.Range("B2").Value = "sdmx-dimension:refArea"
$ws.Range("C2").Value = "iaest-measure:ocupacion-1-digito-descripcion"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("G2").Value = "iaest-measure:sexo"

# Row 3: property kind
$ws.Range("B3").Value = "dim"
$ws.Range("C3").Value = "medida"
$ws.Range("G3").Value = "medida"

# Row 4: datatype / URI template
$ws.Range("B4").Value = "URI-Municipio"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("F4").Value = "URI-Comunidad"
$ws.Range("G4").Value = "xsd:int"

# Row 5 (old per-column mapping-file references) is no longer needed.
$ws.Rows.Item(5).Delete()
